# añadido verificacion de workcenter
#
# Rows 1-2: new H-code + "workcenter" status, with empty C/D/E verification
# columns appended (previously absent).
# Row 3: new H-code + "correcta_por_defecto" status; its old empty C/D/E
# verification columns are removed entirely.
# Rows 4-10: cleared out completely (H-codes/status removed), now matching
# the blank-row pattern used from row 11 onward.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 ---------------------------------------------------------------
$ws.Range("A1").Value = "H59388249"
$ws.Range("B1").Value = "workcenter"
$ws.Range("C1:E1").Value = "'"
$ws.Range("C1:E1").ClearFormats()

# --- Row 2 ---------------------------------------------------------------
$ws.Range("A2").Value = "H59384552"
$ws.Range("B2").Value = "workcenter"
$ws.Range("C2:E2").Value = "'"
$ws.Range("C2:E2").ClearFormats()

# --- Row 3 ---------------------------------------------------------------
$ws.Range("A3").Value = "H59392720"
$ws.Range("B3").Value = "correcta_por_defecto"
$ws.Range("C3:E3").ClearContents()

# --- Rows 4-10: wipe down to the same blank pattern as row 11+ -----------
$ws.Range("A4:E10").Value = "'"
$ws.Range("A4:E10").ClearFormats()
